$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3..53 down to 4..54),
# carrying formatting down from the row above (as Excel does by default).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's new record.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44643
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100112031
$ws.Cells.Item(3, 7).Value = "Poroto verde"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 1700
$ws.Cells.Item(3, 11).Value = 500
$ws.Cells.Item(3, 12).Value = 600
$ws.Cells.Item(3, 13).Value = 550
$ws.Cells.Item(3, 14).Value = "$/kilo"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 550
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
